$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (ID_db_Id_statico_entry) values for rows 35-65 were regenerated
# by the Terminologia_glossario mechanism and are each one less than before.
for ($row = 35; $row -le 65; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value = $cell.Value2 - 1
}
